# Fix R4 to R7 PN
#
# Row 11 of the BOM (designator "R4, R5, R6, R7", comment "2.2k") had the
# wrong manufacturer/part-number pair copied from the R1 row
# (Yageo / RC0603FR-07270RL). The correct part is the Vishay
# CRCW06032K20FKEA (matching the OCTOPART_URL already on that row), so fix
# the MFN (column E) and MPN (column F) cells.
#
# The leading "'" forces Excel's text/quote-prefix so the cell keeps being
# stored as text with the same cell style it already had (these BOM rows
# use quotePrefix="1" styling), instead of picking up a brand-new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "'Vishay"
$ws.Range("F11").Value = "'CRCW06032K20FKEA"
